$d = $word.ActiveDocument
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# ---------------------------------------------------------------------------
# 1. Row "7.76" - split the two notes paragraphs so the spell-checker
#    proofing marks ("dat", "config_parse", "dlgrep") surface as separate
#    runs wrapped in <w:proofErr> start/end markers (same visible text).
# ---------------------------------------------------------------------------
$p163 = $d.Paragraphs(7)
$xml163 = '<w:p ' + $wNs + '>' `
    + '<w:r><w:t>Common IT 163 \u2013 removes RCS revision tags from .</w:t></w:r>' `
    + '<w:proofErr w:type="spellStart"/>' `
    + '<w:r><w:t>dat</w:t></w:r>' `
    + '<w:proofErr w:type="spellEnd"/>' `
    + '<w:r><w:t xml:space="preserve"> files generated by the </w:t></w:r>' `
    + '<w:proofErr w:type="spellStart"/>' `
    + '<w:r><w:t>config_parse</w:t></w:r>' `
    + '<w:proofErr w:type="spellEnd"/>' `
    + '<w:r><w:t xml:space="preserve"> tool</w:t></w:r>' `
    + '</w:p>'
$xml163 = $xml163.Replace('\u2013', [char]0x2013)
$p163.Range.InsertXML($xml163)

$p162 = $d.Paragraphs(8)
$xml162 = '<w:p ' + $wNs + '>' `
    + '<w:r><w:t xml:space="preserve">Common IT 162 \u2013 updates to </w:t></w:r>' `
    + '<w:proofErr w:type="spellStart"/>' `
    + '<w:r><w:t>dlgrep</w:t></w:r>' `
    + '<w:proofErr w:type="spellEnd"/>' `
    + '<w:r><w:t xml:space="preserve"> utility</w:t></w:r>' `
    + '</w:p>'
$xml162 = $xml162.Replace('\u2013', [char]0x2013)
$p162.Range.InsertXML($xml162)

# ---------------------------------------------------------------------------
# 2. Row "7.77" - same proofing-mark treatment around "CGUIButton"; the
#    leading "Common IT 161 - extends" text is unchanged.
# ---------------------------------------------------------------------------
$p161 = $d.Paragraphs(11)
$xml161 = '<w:p ' + $wNs + '>' `
    + '<w:r><w:t>Common IT 161 \u2013 extend</w:t></w:r>' `
    + '<w:r><w:t>s</w:t></w:r>' `
    + '<w:r><w:t xml:space="preserve"> </w:t></w:r>' `
    + '<w:proofErr w:type="spellStart"/>' `
    + '<w:r><w:t>CGUIButton</w:t></w:r>' `
    + '<w:proofErr w:type="spellEnd"/>' `
    + '<w:r><w:t xml:space="preserve"> class to add the state &quot;pressed/disabled&quot;</w:t></w:r>' `
    + '</w:p>'
$xml161 = $xml161.Replace('\u2013', [char]0x2013)
$p161.Range.InsertXML($xml161)

# ---------------------------------------------------------------------------
# 3. Append a new row ("7.78") documenting Common IT 168.
# ---------------------------------------------------------------------------
$table = $d.Tables.Item(1)
$newRow = $table.Rows.Add()
$newRow.Cells.Item(1).Range.Text = "7.78"
# Seed the notes cell with placeholder text first: InsertXML on a still-empty
# paragraph range inserts a sibling paragraph instead of replacing it, so we
# give it real text (and therefore a non-collapsed range) to target.
$newRow.Cells.Item(2).Range.Text = "placeholder"

# The row/cell handles can go stale right after the structural edit above,
# so re-resolve the table/row/cell from the document before touching it again.
$table = $d.Tables.Item(1)
$newRow = $table.Rows.Item($table.Rows.Count)
$cell2 = $newRow.Cells.Item(2)
$p168 = $cell2.Range.Paragraphs.Item(1)
$xml168 = '<w:p ' + $wNs + '>' `
    + '<w:r><w:t xml:space="preserve">Common IT 168 \u2013 correct line number reporting for warnings and errors from </w:t></w:r>' `
    + '<w:proofErr w:type="spellStart"/>' `
    + '<w:r><w:t>config_parse</w:t></w:r>' `
    + '<w:proofErr w:type="spellEnd"/>' `
    + '<w:r><w:t xml:space="preserve"> tool </w:t></w:r>' `
    + '</w:p>'
$xml168 = $xml168.Replace('\u2013', [char]0x2013)
$p168.Range.InsertXML($xml168)
